# Sentinel-2 NG workbook update:
#  - rename Sheet1 -> lplanet2022
#  - add a new sheet vhroda2024 with band-pass info from the VHRODA 2024 presentation
#  - add a source hyperlink on the new sheet
#  - tidy up the selection on the first sheet

$wb = $excel.ActiveWorkbook

# --- rename the existing sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "lplanet2022"

# tidy up the selection on the first sheet before switching away from it
$ws1.Range("D2").Select()

# --- add the new sheet right after it ------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "vhroda2024"

# --- populate the new sheet, matching the exact order the data was typed in ---
# header row: units row first for the numeric columns that got new strings
$ws2.Range("F2").Value = "Lref"
$ws2.Range("G2").Value = "SNR"

# band-name column, in the order the author typed them (skips row 8 initially)
$ws2.Range("A4").Value  = "1(H)"
$ws2.Range("A6").Value  = "2(H)"
$ws2.Range("A11").Value = "4(H)"
$ws2.Range("A12").Value = "5(H)"
$ws2.Range("A13").Value = "6(H)"
$ws2.Range("A14").Value = "7(H)"
$ws2.Range("A15").Value = "8(H)"
$ws2.Range("A16").Value = "8a(H)"
$ws2.Range("A17").Value = "9(H)"
$ws2.Range("A19").Value = "10(H)"
$ws2.Range("A20").Value = "11(H)"
$ws2.Range("A8").Value  = "3(H)"

# remaining header cells
$ws2.Range("B2").Value = "current[m]"
$ws2.Range("C2").Value = "SSD[m]"
$ws2.Range("D2").Value = "CWvl[nm]"
$ws2.Range("E2").Value = "BW[nm]"

# the reference / source link, entered last
$ws2.Range("A1").Value = "https://earth.esa.int/eogateway/documents/d/earth-online/15_sproud_esa"
$ws2.Hyperlinks.Add($ws2.Range("A1"), "https://earth.esa.int/eogateway/documents/d/earth-online/15_sproud_esa")

# band names that re-use already existing shared strings
$ws2.Range("A2").Value  = "Band"
$ws2.Range("A5").Value  = "1a"
$ws2.Range("A7").Value  = "2a"
$ws2.Range("A9").Value  = "3a"
$ws2.Range("A10").Value = "3b"
$ws2.Range("A18").Value = "9a"
$ws2.Range("A21").Value = "12a"
$ws2.Range("A22").Value = "12b"
$ws2.Range("A23").Value = "12c"

# numeric data
$ws2.Range("A3").Value = 0

$ws2.Range("C3").Value = 20
$ws2.Range("D3").Value = 412
$ws2.Range("E3").Value = 20
$ws2.Range("F3").Value = 45
$ws2.Range("G3").Value = 130

$ws2.Range("B4").Value = 60
$ws2.Range("C4").Value = 20
$ws2.Range("D4").Value = 443
$ws2.Range("E4").Value = 20
$ws2.Range("F4").Value = 129.11
$ws2.Range("G4").Value = 332

$ws2.Range("C5").Value = 20
$ws2.Range("D5").Value = 475
$ws2.Range("E5").Value = 20
$ws2.Range("F5").Value = 90.7
$ws2.Range("G5").Value = 190

$ws2.Range("B6").Value = 10
$ws2.Range("C6").Value = 5
$ws2.Range("D6").Value = 490
$ws2.Range("E6").Value = 65
$ws2.Range("F6").Value = 128
$ws2.Range("G6").Value = 106

$ws2.Range("C7").Value = 10
$ws2.Range("D7").Value = 520
$ws2.Range("E7").Value = 15
$ws2.Range("F7").Value = 80.5
$ws2.Range("G7").Value = 120

$ws2.Range("B8").Value = 10
$ws2.Range("C8").Value = 5
$ws2.Range("D8").Value = 560
$ws2.Range("E8").Value = 35
$ws2.Range("F8").Value = 128
$ws2.Range("G8").Value = 121

$ws2.Range("C9").Value = 20
$ws2.Range("D9").Value = 620
$ws2.Range("E9").Value = 30
$ws2.Range("F9").Value = 29.8
$ws2.Range("G9").Value = 200

$ws2.Range("C10").Value = 20
$ws2.Range("D10").Value = 650
$ws2.Range("E10").Value = 20
$ws2.Range("F10").Value = 29.8
$ws2.Range("G10").Value = 200

$ws2.Range("B11").Value = 10
$ws2.Range("C11").Value = 5
$ws2.Range("D11").Value = 665
$ws2.Range("E11").Value = 30
$ws2.Range("F11").Value = 108
$ws2.Range("G11").Value = 115

$ws2.Range("B12").Value = 20
$ws2.Range("C12").Value = 10
$ws2.Range("D12").Value = 705
$ws2.Range("E12").Value = 15
$ws2.Range("F12").Value = 74.6
$ws2.Range("G12").Value = 124

$ws2.Range("B13").Value = 20
$ws2.Range("C13").Value = 10
$ws2.Range("D13").Value = 740
$ws2.Range("E13").Value = 15
$ws2.Range("F13").Value = 68.23
$ws2.Range("G13").Value = 111

$ws2.Range("B14").Value = 20
$ws2.Range("C14").Value = 10
$ws2.Range("D14").Value = 783
$ws2.Range("E14").Value = 20
$ws2.Range("F14").Value = 66.7
$ws2.Range("G14").Value = 116

$ws2.Range("B15").Value = 10
$ws2.Range("C15").Value = 5
$ws2.Range("D15").Value = 842
$ws2.Range("E15").Value = 115
$ws2.Range("F15").Value = 103
$ws2.Range("G15").Value = 115

$ws2.Range("B16").Value = 20
$ws2.Range("C16").Value = 10
$ws2.Range("D16").Value = 865
$ws2.Range("E16").Value = 20
$ws2.Range("F16").Value = 52.39
$ws2.Range("G16").Value = 84

$ws2.Range("B17").Value = 60
$ws2.Range("C17").Value = 60
$ws2.Range("D17").Value = 945
$ws2.Range("E17").Value = 20
$ws2.Range("F17").Value = 8.77
$ws2.Range("G17").Value = 166

$ws2.Range("C18").Value = 60
$ws2.Range("D18").Value = 985
$ws2.Range("E18").Value = 20
$ws2.Range("F18").Value = 56.4
$ws2.Range("G18").Value = 90

$ws2.Range("B19").Value = 60
$ws2.Range("C19").Value = 60
$ws2.Range("D19").Value = 1375
$ws2.Range("E19").Value = 30
$ws2.Range("F19").Value = 6
$ws2.Range("G19").Value = 400

$ws2.Range("B20").Value = 20
$ws2.Range("C20").Value = 10
$ws2.Range("D20").Value = 1610
$ws2.Range("E20").Value = 90
$ws2.Range("F20").Value = 4
$ws2.Range("G20").Value = 84

$ws2.Range("B21").Value = 20
$ws2.Range("C21").Value = 10
$ws2.Range("D21").Value = 2130
$ws2.Range("E21").Value = 50
$ws2.Range("F21").Value = 1.7
$ws2.Range("G21").Value = 50

$ws2.Range("B22").Value = 20
$ws2.Range("C22").Value = 10
$ws2.Range("D22").Value = 2210
$ws2.Range("E22").Value = 50
$ws2.Range("F22").Value = 1.7
$ws2.Range("G22").Value = 50

$ws2.Range("B23").Value = 20
$ws2.Range("C23").Value = 10
$ws2.Range("D23").Value = 2260
$ws2.Range("E23").Value = 50
$ws2.Range("F23").Value = 1.7
$ws2.Range("G23").Value = 50

# selection / active cell on the new sheet (also leaves vhroda2024 as the active tab)
$ws2.Range("L7").Select()
